$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 holds the "FolderExpName" header values which are the Experiment
# result-folder paths. Columns B:Y share one string (Experiment_A1/...)
# and columns Z:AW share another (Experiment_A2/...). Update the text in
# place (same cells / same shared-string slots) so the generated
# OpenScenario result folders get the "_generated" suffix.
$oldA1 = "KTH_pedestrian_autoware_light/OpenScenario/Results/Experiment_A1/OpenScenario"
$newA1 = "KTH_pedestrian_autoware_light/OpenScenario/Results/Experiment_A1_generated/OpenScenario"
$oldA2 = "KTH_pedestrian_autoware_light/OpenScenario/Results/Experiment_A2/OpenScenario"
$newA2 = "KTH_pedestrian_autoware_light/OpenScenario/Results/Experiment_A2_generated/OpenScenario"

for ($col = 2; $col -le 25; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    if ($cell.Value2 -eq $oldA1) {
        $cell.Value2 = $newA1
    }
}

for ($col = 26; $col -le 49; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    if ($cell.Value2 -eq $oldA2) {
        $cell.Value2 = $newA2
    }
}

# The longer text no longer fits the previous best-fit column width, so
# the columns widen accordingly (closest attainable width to 82.109375).
$ws.Range("B1:AW8").EntireColumn.ColumnWidth = 81.26
